$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

for ($r = 1; $r -le 319; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq $oldVal) {
        $cell.Value = $newVal
    }
}
